# cy_hist.xlsx — "some improvements to image floating / updated sizes in
# histogram spreadsheets": the histogram chart on the sheet is being resized
# (made smaller) in place; its top-left anchor stays put, only its
# bottom-right corner moves.
#
# Diff says the <xdr:to> anchor goes from
#     col=12 colOff=313649  row=21 rowOff=113850
# to
#     col=9  colOff=342449  row=14 rowOff=7350
# while <xdr:from> (col 3 / colOff 400049 / row 2 / rowOff 133350) is
# untouched. Driving this through the ChartObject's Left/Top/Width/Height
# (rather than poking the drawing XML directly) is the COM-idiomatic way to
# "resize a chart" and lets the host recompute the two-cell anchor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$co = $ws.ChartObjects().Item(1)

# Left/Top (driven by <xdr:from>) are unchanged by this edit -- only the
# size shrinks, which moves <xdr:to>.
$co.Left = 273.661054072342
$co.Top  = 40.5
$co.Width  = 346.08956692913387
$co.Height = 170.07874015748033

# Cosmetic-only: the workbook's saved window position also moved a little
# (yWindow 1800 -> 2400) in the source commit. Reflect the intent on the
# live window even though this is purely a UI/view setting.
try {
    $excel.ActiveWindow.Top = 2400
} catch {
}
